# "fix excel 3jr/H=1 credit"
# The "Crédit" column of Tableau3 was computed as
#   Charge de travail (en jour/Homme) / 2
# but per the NOTE on the sheet, 3 jours/Homme = 1 crédit, so every
# formula in that calculated column (and its calculatedColumnFormula
# definition) must divide by 3 instead of 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update every data cell of the "Crédit" column (F2:F19) individually so
# each keeps its own full structured-reference formula (matching how the
# table's calculated column is stored) rather than collapsing into a
# single shared formula.
for ($r = 2; $r -le 19; $r++) {
    $ws.Range("F$r").Formula = "=[@[Charge de travail (en jour/Homme)]]/3"
}

# Restore the zoom level and move the selection to F6, reflecting the
# view state captured the last time the workbook was saved.
$excel.ActiveWindow.Zoom = 100
$ws.Range("F6").Select()
